$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "301.17"
Set-TextValue $ws.Range("E2") "0.19%"
Set-TextValue $ws.Range("D3") "31.75"
Set-TextValue $ws.Range("E3") "1.15%"
Set-TextValue $ws.Range("D4") "5.099"
Set-TextValue $ws.Range("E4") "-0.90%"
Set-TextValue $ws.Range("D5") "0.08199"
Set-TextValue $ws.Range("E5") "11.12%"
Set-TextValue $ws.Range("D6") "2.608"
Set-TextValue $ws.Range("E6") "4.22%"
Set-TextValue $ws.Range("D7") "7.827"
Set-TextValue $ws.Range("E7") "-1.37%"
Set-TextValue $ws.Range("D8") "3.842"
Set-TextValue $ws.Range("E8") "2.01%"
Set-TextValue $ws.Range("D9") "0.9249"
Set-TextValue $ws.Range("E9") "0.26%"
Set-TextValue $ws.Range("D10") "0.1753"
Set-TextValue $ws.Range("E10") "1.37%"
Set-TextValue $ws.Range("D11") "0.07479"
Set-TextValue $ws.Range("E11") "-1.44%"
Set-TextValue $ws.Range("D12") "0.08904"
Set-TextValue $ws.Range("E12") "9.72%"
Set-TextValue $ws.Range("D13") "0.02997"
Set-TextValue $ws.Range("E13") "-1.14%"
Set-TextValue $ws.Range("D14") "0.09999"
Set-TextValue $ws.Range("E14") "0.74%"
Set-TextValue $ws.Range("D15") "0.001506"
Set-TextValue $ws.Range("E15") "0.90%"
Set-TextValue $ws.Range("D16") "0.005794"
Set-TextValue $ws.Range("E16") "-4.74%"
Set-TextValue $ws.Range("D17") "3.592"
Set-TextValue $ws.Range("E17") "3.82%"
Set-TextValue $ws.Range("D18") "2.261"
Set-TextValue $ws.Range("E18") "1.43%"
Set-TextValue $ws.Range("E19") "-1.03%"
Set-TextValue $ws.Range("D20") "0.1343"
Set-TextValue $ws.Range("E20") "0.51%"
Set-TextValue $ws.Range("D21") "3.907"
Set-TextValue $ws.Range("E21") "-16.09%"
Set-TextValue $ws.Range("D22") "0.1696"
Set-TextValue $ws.Range("E22") "8.30%"
Set-TextValue $ws.Range("E23") "-0.97%"
Set-TextValue $ws.Range("D24") "0.001243"
Set-TextValue $ws.Range("E24") "1.51%"
Set-TextValue $ws.Range("D25") "0.004527"
Set-TextValue $ws.Range("E25") "0.76%"
Set-TextValue $ws.Range("D26") "0.0001198"
Set-TextValue $ws.Range("E26") "-7.84%"
Set-TextValue $ws.Range("D27") "0.0003403"
Set-TextValue $ws.Range("E27") "81.81%"
Set-TextValue $ws.Range("D39") "0.01768"
Set-TextValue $ws.Range("E39") "2.16%"
Set-TextValue $ws.Range("D40") "0.04541"
Set-TextValue $ws.Range("E40") "0.32%"
Set-TextValue $ws.Range("D41") "0.006979"
Set-TextValue $ws.Range("E41") "-2.91%"
Set-TextValue $ws.Range("D42") "0.1377"
Set-TextValue $ws.Range("E42") "2.19%"
Set-TextValue $ws.Range("D43") "0.002137"
Set-TextValue $ws.Range("E43") "-4.12%"
Set-TextValue $ws.Range("D44") "0.009608"
Set-TextValue $ws.Range("E44") "-10.46%"
Set-TextValue $ws.Range("D45") "0.00006386"
Set-TextValue $ws.Range("E45") "1.71%"
Set-TextValue $ws.Range("D46") "0.00000000748"
Set-TextValue $ws.Range("E46") "-0.22%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue $ws.Range("D47") "0.8040"
Set-TextValue $ws.Range("E47") "-0.53%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue $ws.Range("D48") "0.008739"
Set-TextValue $ws.Range("E48") "-12.63%"
Set-TextValue $ws.Range("D49") "0.00002096"
Set-TextValue $ws.Range("E49") "-0.22%"
Set-TextValue $ws.Range("D50") "0.0001996"
Set-TextValue $ws.Range("E50") "-0.15%"
